# Insert a new data row above row 572 (pushes existing rows 572..627 down to
# 573..628) and populate it with a new weekly price record for
# "Feria Lagunitas de Puerto Montt" / Zanahoria.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 572; this shifts all rows below
# down by one and grows the sheet's used range from R627 to R628.
$ws.Rows.Item(572).Insert()

# Populate the newly inserted row 572 with the new record's values.
$ws.Cells.Item(572, 1).Value  = 4
$ws.Cells.Item(572, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(572, 3).Value  = "Los Lagos"
$ws.Cells.Item(572, 4).Value  = 45166
$ws.Cells.Item(572, 5).Value  = 10
$ws.Cells.Item(572, 6).Value  = 100114013
$ws.Cells.Item(572, 7).Value  = "Zanahoria"
$ws.Cells.Item(572, 8).Value  = "Sin especificar"
$ws.Cells.Item(572, 9).Value  = "Primera"
$ws.Cells.Item(572, 10).Value = 250
$ws.Cells.Item(572, 11).Value = 7500
$ws.Cells.Item(572, 12).Value = 9000
$ws.Cells.Item(572, 13).Value = 8100
$ws.Cells.Item(572, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(572, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(572, 16).Value = 405
$ws.Cells.Item(572, 17).Value = 20
$ws.Cells.Item(572, 18).Value = "Hortaliza"
